$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (index 3) for rows 2..409: 45203 -> 45204 ("Förändrad" date)
for ($i = 2; $i -le 409; $i++) {
    $cell = $ws.Cells.Item($i, 3)
    $cell.Value2 = 45204
}

# Row 409 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(409).RowHeight = 15

# Append new row 410 with the new logging-notice record
$ws.Cells.Item(410, 1).Value = "A 47593-2023"

$ws.Cells.Item(410, 2).Value2 = 45203
$ws.Cells.Item(410, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(410, 3).Value2 = 45204
$ws.Cells.Item(410, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(410, 4).Value = "VÄSTERBOTTENS LÄN"
$ws.Cells.Item(410, 5).Value = "MALÅ"
$ws.Cells.Item(410, 6).Value = "Sveaskog"

$ws.Cells.Item(410, 7).Value = 23.1
$ws.Cells.Item(410, 8).Value = 0
$ws.Cells.Item(410, 9).Value = 0
$ws.Cells.Item(410, 10).Value = 0
$ws.Cells.Item(410, 11).Value = 0
$ws.Cells.Item(410, 12).Value = 0
$ws.Cells.Item(410, 13).Value = 0
$ws.Cells.Item(410, 14).Value = 0
$ws.Cells.Item(410, 15).Value = 0
$ws.Cells.Item(410, 16).Value = 0
$ws.Cells.Item(410, 17).Value = 0

# R410 mirrors the blank, wrap-text styled "Artnamn" cells used throughout the sheet
$ws.Cells.Item(410, 18).WrapText = $true

Write-Host "done"
